# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# 1. Update "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 7 de Mayo de 2020 a las 01:04"

# 2. Estados Unidos (row 4)
$ws.Range("B4").Value = 1257156
$ws.Range("C4").Value = 19523
$ws.Range("E4").Value = 977746
$ws.Range("G4").Value = 1871
$ws.Range("H4").Value = 74142

# 3. Chequia (row 48)
$ws.Range("B48").Value = 7974
$ws.Range("C48").Value = 78
$ws.Range("D48").Value = 4205
$ws.Range("E48").Value = 3507

# 4. Republica de Yibuti (row 91)
$ws.Range("B91").Value = 1124
$ws.Range("C91").Value = 4
$ws.Range("D91").Value = 755
$ws.Range("E91").Value = 366
$ws.Range("G91").Value = 1
$ws.Range("H91").Value = 3

# 5. Swap Seychelles (row 205) and Montserrat (row 206): they exchange rank
#    position, so the whole rows (name + all stats) trade places, then the
#    new Montserrat (row 205) and Seychelles (row 206) rows get refreshed
#    stats too.
$ws.Range("A205").Value = "Montserrat"
$ws.Range("B205").Value = 11
$ws.Range("C205").Value = 0
$ws.Range("D205").Value = 7
$ws.Range("E205").Value = 3
$ws.Range("F205").Value = 1
$ws.Range("G205").Value = 0
$ws.Range("H205").Value = 1

$ws.Range("A206").Value = "Seychelles"
$ws.Range("B206").Value = 11
$ws.Range("C206").Value = 0
$ws.Range("D206").Value = 8
$ws.Range("E206").Value = 3
$ws.Range("F206").Value = 0
$ws.Range("G206").Value = 0
$ws.Range("H206").Value = 0
